$d = $word.ActiveDocument
$CR = [char]13

# --- Locate the three "FECOLn : fecolN" paragraphs and append " esquema fecolN" ---
$fecol3Index = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd($CR)
    if ($t -eq "FECOL1 : fecol1") {
        $p.Range.InsertAfter(" esquema fecol1")
    } elseif ($t -eq "FECOL2 : fecol2") {
        $p.Range.InsertAfter(" esquema fecol2")
    } elseif ($t -eq "FECOL3 : fecol3") {
        $p.Range.InsertAfter(" esquema fecol3")
        $fecol3Index = $i
    }
}

# --- Insert a brand new paragraph right after the FECOL3 paragraph ---
$p3 = $d.Paragraphs($fecol3Index)
$p3.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($fecol3Index + 1)
$fullNew = $d.Range($newPara.Range.Start, $newPara.Range.End)
$uAcc = [char]250
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:r><w:t xml:space="preserve">Cada usuario tiene asignado un esquema donde puede crear objetos </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>sql</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>. Por favor no crear ning' + $uAcc + 'n objeto en el esquema dbo.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
$fullNew.InsertXML($newParaXml)

# --- Drop the stray <w:lastRenderedPageBreak/> before "CfdiColombiaSqlSetup" ---
# Re-writing the matched text via Find/Replace causes the run to be regenerated
# without the obsolete lastRenderedPageBreak marker (proofErr wrappers untouched).
$d.Content.Find.Execute("CfdiColombiaSqlSetup", $true, $false, $false, $false, $false,
                         $true, 1, $false, "CfdiColombiaSqlSetup", 2)

# --- Turn the trailing _GoBack bookmark paragraph back into a plain empty paragraph ---
$lastIdx = $d.Paragraphs.Count
$lastP = $d.Paragraphs($lastIdx)
$fullLast = $d.Range($lastP.Range.Start, $lastP.Range.End)
$fullLast.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>')
